# The source data ("Índice do emprego formal") was refreshed with a new
# year (2023) added for every region (Brasil, Nordeste, Sergipe). The
# underlying table is grouped by region and sorted by date within each
# region, so inserting the new 2023 rows pushes every following region's
# block down by one row; the sheet grows from 31 to 34 data+header rows
# (A1:D31 -> A1:D34).
#
# Rather than performing a sequence of row inserts/shifts, we simply
# rewrite the whole data block (rows 2-34) with the final values taken
# from the updated table - this reproduces exactly the same end state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2013', 100),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2014', 101.2729253253112),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2015', 99.86379134956168),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2016', 94.09943317286582),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2017', 94.55172957222145),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2018', 95.26579737496398),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2019', 103.3813605432476),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2020', 100.5542731061483),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2021', 99.55144222900864),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2022', 107.8499571171155),
    @('Brasil', 'Índice do emprego formal: 2013=100', '31/12/2023', 111.7633020039681),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2013', 100),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2014', 102.3093950626827),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2015', 99.69270873591725),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2016', 94.50517603910063),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2017', 95.70884457991802),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2018', 96.86924970117771),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2019', 95.76212288737956),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2020', 93.74482872189195),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2021', 101.1677314486524),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2022', 109.5253234394307),
    @('Nordeste', 'Índice do emprego formal: 2013=100', '31/12/2023', 115.6776572779893),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2013', 100),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2014', 102.7719795453145),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2015', 99.80112131107141),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2016', 94.40527385866552),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2017', 96.22623375023103),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2018', 95.95243669521287),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2019', 86.96198632246936),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2020', 90.21107756761754),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2021', 96.62251247612595),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2022', 102.7995810486107),
    @('Sergipe', 'Índice do emprego formal: 2013=100', '31/12/2023', 111.417411126856)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
